# Lab01_ReviewReport.xlsx edit:
#  - Sheet "Requirements Phase Defects": fill in E12 ("Variable types was not
#    stated") review comment and move the selection there.
#  - Sheet "Architect. Design Phase Defects": no longer the active tab.
#  - Sheet "Coding Phase Defects": becomes the active tab, gets review
#    comments filled in for rows 10-16 (columns C/D/E) and selection moves
#    to E16.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Requirements Phase Defects ------------------------------
$wsReq = $wb.Worksheets.Item("Requirements Phase Defects")
$wsReq.Range("E12").Value = "Variable types was not stated"
$wsReq.Rows.Item(12).RowHeight = 13.8
# Selection is set now (sheet still inactive); activating another sheet
# later does not disturb a sheet's own remembered selection, but
# selecting on this sheet AFTER activating a different one would steal
# activation back - so do it before switching tabs below.
[void]$wsReq.Range("E12").Select()

# --- Sheet 3: Coding Phase Defects -------------------------------------
$wsCode = $wb.Worksheets.Item("Coding Phase Defects")

$wsCode.Range("C10").Value = "C01"
$wsCode.Range("D10").Value = "CartiRepoMock line 80"
$wsCode.Range("E10").Value = "Se presupune sa ia cartile dintr-un anumit an da le ia alea care nu is in anu ala"

$wsCode.Range("C11").Value = "C01"
$wsCode.Range("D11").Value = "Validator line 16-19"
$wsCode.Range("E11").Value = "Aci tre sa verifice ca lista ii vida, nu nula"

$wsCode.Range("C12").Value = "C02"
$wsCode.Range("D12").Value = "CartiRepoMock line 80"
$wsCode.Range("E12").Value = "Se presupune sa ia cartile dintr-un anumit an da le ia alea care nu is in anu ala"

$wsCode.Range("C13").Value = "C06"
$wsCode.Range("D13").Value = "some places"
$wsCode.Range("E13").Value = "There is checking with regex [10-9]+ which is equivalent with [0-9]+ where I presume the desired check would be with 1[0-9]+ (those checks are in the year)"
$wsCode.Rows.Item(13).RowHeight = 46.25

$wsCode.Range("C14").Value = "C09"
$wsCode.Range("D14").Value = "some places"
$wsCode.Range("E14").Value = "errors in canfusing function names isOKString and isStringOK, in some places authors were named referents, functiile cautaDupa* din clasa Carte is confuze rau de tot"
$wsCode.Rows.Item(14).RowHeight = 46.25

$wsCode.Range("C15").Value = "C10"
$wsCode.Range("D15").Value = "in repo"
$wsCode.Range("E15").Value = "while i<= size should’ve been I < size"
$wsCode.Rows.Item(15).RowHeight = 13.8

$wsCode.Range("C16").Value = "C12"
$wsCode.Range("D16").Value = "in Carte"
$wsCode.Range("E16").Value = "an could’ve been integer not string"
$wsCode.Rows.Item(16).RowHeight = 13.8

# Coding Phase Defects becomes the active / visible tab, with its
# selection on E16. Do this last, since selecting/activating a sheet
# makes it the active one.
[void]$wsCode.Activate()
[void]$wsCode.Range("E16").Select()
